# 835 Medium Bit OP Class:?  -- add " Image Overlap" as a new tracked row (26)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bring row 26 up to the same formatting as the other data rows -------
# Row 21 has the exact same cell layout we need (A:E then G:J, no F, no K),
# so copy its formats first and then overwrite with the new row's values.
$ws.Range("A21:E21").Copy()
$ws.Range("A26:E26").PasteSpecial(-4122)
$ws.Range("G21:J21").Copy()
$ws.Range("G26:J26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- New row 26 content ----------------------------------------------------
$ws.Range("A26").Value = "?"
$ws.Range("B26").Value = " Image Overlap"
$ws.Range("C26").Value = 835
$ws.Range("D26").Value = "Medium"
$ws.Range("E26").Value = "Bit OP"
$ws.Range("G26").Value = 0.86111111111111116
$ws.Range("H26").Value = 0.95833333333333337
$ws.Range("I26").Formula = "=H26-G26"
$ws.Range("J26").Value = "TIME EXCEED"

# --- View state: scrolled down, K29 selected -------------------------------
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K29").Select()
